$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as plain text in this workbook
# (e.g. "60.924.78", "527.60", "0.0₃0871"). Force every Price cell we
# touch to Text format first so Excel does not auto-convert the
# numeric-looking ones to a Double (which would silently drop
# significant trailing zeros, e.g. "527.60" -> 527.6).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '60.924.78'
$ws.Range('E2').Value = '  -1.49%  '

$ws.Range('D3').Value = '2.911.11'
$ws.Range('E3').Value = '  -2.70%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '527.60'
$ws.Range('E5').Value = '  -2.63%  '

$ws.Range('D6').Value = '144.70'
$ws.Range('E6').Value = '  -5.10%  '

$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').Value = '0.548'
$ws.Range('E8').Value = '  -3.56%  '

$ws.Range('D9').Value = '2.918.49'
$ws.Range('E9').Value = '  -2.90%  '

$ws.Range('D10').Value = '0.108'
$ws.Range('E10').Value = '  -4.49%  '

$ws.Range('D11').Value = '6.15'
$ws.Range('E11').Value = '  -0.23%  '

$ws.Range('D12').Value = '0.359'
$ws.Range('E12').Value = '  -2.60%  '

$ws.Range('D13').Value = '3.418.55'
$ws.Range('E13').Value = '  -2.71%  '

$ws.Range('E14').Value = '  +2.89%  '

$ws.Range('D15').Value = '60.855.40'
$ws.Range('E15').Value = '  -1.69%  '

$ws.Range('D16').Value = '22.56'
$ws.Range('E16').Value = '  -5.58%  '

$ws.Range('D17').Value = '2.904.61'
$ws.Range('E17').Value = '  -3.12%  '

$ws.Range('D18').Value = '0.0000142'
$ws.Range('E18').Value = '  -3.43%  '

$ws.Range('E19').Value = '  -5.13%  '

$ws.Range('D20').Value = '11.61'
$ws.Range('E20').Value = '  -3.74%  '

$ws.Range('D21').Value = '354.55'
$ws.Range('E21').Value = '  -6.58%  '

$ws.Range('D22').Value = '6.53'
$ws.Range('E22').Value = '  -3.45%  '

$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.22%  '

$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Value = '64.95'
$ws.Range('E24').Value = '  -1.57%  '

$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').Value = '0.452'
$ws.Range('E25').Value = '  -3.82%  '

$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.177'
$ws.Range('E26').Value = '  -6.38%  '

$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '7.88'
$ws.Range('E28').Value = '  -4.37%  '

$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0871'
$ws.Range('E29').Value = '  -7.41%  '

$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  +0.02%  '

$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.69'
$ws.Range('E31').Value = '  -2.33%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '19.69'
$ws.Range('E32').Value = '  -3.76%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '152.97'
$ws.Range('E33').Value = '  -5.10%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '4.42'
$ws.Range('E34').Value = '  -3.99%  '

$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '5.57'
$ws.Range('E35').Value = '  -6.79%  '

$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -7.09%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.20'
$ws.Range('E37').Value = '  -6.42%  '

$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = '37.65'
$ws.Range('E38').Value = '  -0.13%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.47'
$ws.Range('E39').Value = '  -4.75%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '3.71'
$ws.Range('E40').Value = '  -4.98%  '

$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '0.653'
$ws.Range('E41').Value = '  -2.99%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.282.19'
$ws.Range('E42').Value = '  -5.85%  '

$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = '0.0583'
$ws.Range('E43').Value = '  -1.49%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '20.38'
$ws.Range('E44').Value = '  -8.00%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '0.997'
$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '4.94'
$ws.Range('E46').Value = '  -4.39%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0238'
$ws.Range('E47').Value = '  -3.08%  '

$ws.Range('B48').Value = 'WhiteBITCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D48').Value = '10.33'
$ws.Range('E48').Value = '  -0.77%  '

$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.0919'
$ws.Range('E49').Value = '  -3.66%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '18.52'
$ws.Range('E50').Value = '  -6.59%  '

$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').Value = '247.65'
$ws.Range('E51').Value = '  -7.67%  '
